$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 3000.1875
$ws.Cells.Item(64, 10).Value = 3167.1667
$ws.Cells.Item(64, 12).Value = 3167.1667
$ws.Cells.Item(64, 14).Value = -3663.1667
$ws.Cells.Item(67, 8).Value = 3000.1875
$ws.Cells.Item(67, 10).Value = 3167.1667
$ws.Cells.Item(67, 12).Value = 3167.1667
$ws.Cells.Item(67, 14).Value = -4883.1667
$ws.Cells.Item(74, 8).Value = 4180.4287
$ws.Cells.Item(74, 9).Value = 3751.5
$ws.Cells.Item(74, 10).Value = 4352
$ws.Cells.Item(74, 11).Value = 3751.5
$ws.Cells.Item(74, 12).Value = 4352
$ws.Cells.Item(74, 13).Value = -2815.5
$ws.Cells.Item(74, 14).Value = -6224
$ws.Cells.Item(77, 8).Value = 4180.4287
$ws.Cells.Item(77, 9).Value = 3751.5
$ws.Cells.Item(77, 10).Value = 4352
$ws.Cells.Item(77, 11).Value = 18757.5
$ws.Cells.Item(77, 12).Value = 21760
$ws.Cells.Item(77, 13).Value = -14077.5
$ws.Cells.Item(77, 14).Value = -31120
$ws.Cells.Item(80, 8).Value = 611305.8
$ws.Cells.Item(80, 9).Value = 375.42856
$ws.Cells.Item(80, 10).Value = 1038957.1
$ws.Cells.Item(80, 11).Value = 1126.28568
$ws.Cells.Item(80, 12).Value = 3116871.3
$ws.Cells.Item(80, 13).Value = -128.28568
$ws.Cells.Item(80, 14).Value = -3118867.3
$ws.Cells.Item(83, 8).Value = 611305.8
$ws.Cells.Item(83, 9).Value = 375.42856
$ws.Cells.Item(83, 10).Value = 1038957.1
$ws.Cells.Item(83, 11).Value = 3378.85704
$ws.Cells.Item(83, 12).Value = 9350613.9
$ws.Cells.Item(83, 13).Value = 1613.14296
$ws.Cells.Item(83, 14).Value = -9360597.9
$ws.Cells.Item(86, 8).Value = 22060.6
$ws.Cells.Item(86, 9).Value = 22060.6
$ws.Cells.Item(86, 11).Value = 22060.6
$ws.Cells.Item(86, 13).Value = -20937.6
$ws.Cells.Item(89, 8).Value = 22060.6
$ws.Cells.Item(89, 9).Value = 22060.6
$ws.Cells.Item(89, 11).Value = 110303
$ws.Cells.Item(89, 13).Value = -104687

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 29415342
$ws.Cells.Item(132, 9).Value = 35717508
$ws.Cells.Item(132, 10).Value = 5233
$ws.Cells.Item(132, 11).Value = 107152524
$ws.Cells.Item(132, 12).Value = 15699
$ws.Cells.Item(132, 13).Value = -107149994
$ws.Cells.Item(132, 14).Value = -20759

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 1371.4242
$ws.Cells.Item(105, 9).Value = 1160.32
$ws.Cells.Item(105, 10).Value = 2031.125
$ws.Cells.Item(105, 11).Value = 1160.32
$ws.Cells.Item(105, 12).Value = 2031.125
$ws.Cells.Item(105, 13).Value = 586.6800000000001
$ws.Cells.Item(105, 14).Value = -5525.125

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 66
$ws.Cells.Item(7, 9).Value = 32.57143
$ws.Cells.Item(7, 10).Value = 300
$ws.Cells.Item(7, 11).Value = 32.57143
$ws.Cells.Item(7, 12).Value = 300
$ws.Cells.Item(7, 13).Value = 80.42857000000001
$ws.Cells.Item(7, 14).Value = -526
$ws.Cells.Item(31, 8).Value = 2492.4412
$ws.Cells.Item(31, 9).Value = 1831
$ws.Cells.Item(31, 10).Value = 3153.8823
$ws.Cells.Item(31, 11).Value = 1831
$ws.Cells.Item(31, 12).Value = 3153.8823
$ws.Cells.Item(31, 13).Value = -1536
$ws.Cells.Item(31, 14).Value = -3743.8823
$ws.Cells.Item(34, 8).Value = 2492.4412
$ws.Cells.Item(34, 9).Value = 1831
$ws.Cells.Item(34, 10).Value = 3153.8823
$ws.Cells.Item(34, 11).Value = 1831
$ws.Cells.Item(34, 12).Value = 3153.8823
$ws.Cells.Item(34, 13).Value = -1629
$ws.Cells.Item(34, 14).Value = -3557.8823
$ws.Cells.Item(62, 8).Value = 5097.5
$ws.Cells.Item(62, 9).Value = 2926.6667
$ws.Cells.Item(62, 10).Value = 6400
$ws.Cells.Item(62, 11).Value = 2926.6667
$ws.Cells.Item(62, 12).Value = 6400
$ws.Cells.Item(62, 13).Value = -2302.6667
$ws.Cells.Item(62, 14).Value = -7648
$ws.Cells.Item(65, 8).Value = 5097.5
$ws.Cells.Item(65, 9).Value = 2926.6667
$ws.Cells.Item(65, 10).Value = 6400
$ws.Cells.Item(65, 11).Value = 14633.3335
$ws.Cells.Item(65, 12).Value = 32000
$ws.Cells.Item(65, 13).Value = -11513.3335
$ws.Cells.Item(65, 14).Value = -38240
$ws.Cells.Item(107, 8).Value = 2269.7334
$ws.Cells.Item(107, 9).Value = 623
$ws.Cells.Item(107, 11).Value = 623
$ws.Cells.Item(107, 13).Value = 1297
$ws.Cells.Item(134, 8).Value = 2295.074
$ws.Cells.Item(134, 9).Value = 1441.2222
$ws.Cells.Item(134, 10).Value = 4002.7778
$ws.Cells.Item(134, 11).Value = 4323.6666
$ws.Cells.Item(134, 12).Value = 12008.3334
$ws.Cells.Item(134, 13).Value = -1788.6666
$ws.Cells.Item(134, 14).Value = -17078.3334

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(62, 8).Value = 3799
$ws.Cells.Item(62, 10).Value = 3998.9285
$ws.Cells.Item(62, 12).Value = 11996.7855
$ws.Cells.Item(62, 14).Value = -13368.7855
$ws.Cells.Item(63, 8).Value = 2460.3333
$ws.Cells.Item(63, 9).Value = 1004
$ws.Cells.Item(63, 10).Value = 3916.6667
$ws.Cells.Item(63, 11).Value = 3012
$ws.Cells.Item(63, 12).Value = 11750.0001
$ws.Cells.Item(63, 13).Value = -2263
$ws.Cells.Item(63, 14).Value = -13248.0001
$ws.Cells.Item(65, 8).Value = 3799
$ws.Cells.Item(65, 10).Value = 3998.9285
$ws.Cells.Item(65, 12).Value = 35990.3565
$ws.Cells.Item(65, 14).Value = -42854.3565
$ws.Cells.Item(66, 8).Value = 2460.3333
$ws.Cells.Item(66, 9).Value = 1004
$ws.Cells.Item(66, 10).Value = 3916.6667
$ws.Cells.Item(66, 11).Value = 9036
$ws.Cells.Item(66, 12).Value = 35250.0003
$ws.Cells.Item(66, 13).Value = -5292
$ws.Cells.Item(66, 14).Value = -42738.0003
$ws.Cells.Item(75, 8).Value = 2579.2727
$ws.Cells.Item(75, 10).Value = 3966.6667
$ws.Cells.Item(75, 12).Value = 11900.0001
$ws.Cells.Item(75, 14).Value = -13896.0001
$ws.Cells.Item(78, 8).Value = 2579.2727
$ws.Cells.Item(78, 10).Value = 3966.6667
$ws.Cells.Item(78, 12).Value = 35700.0003
$ws.Cells.Item(78, 14).Value = -45684.0003
$ws.Cells.Item(81, 8).Value = 27300
$ws.Cells.Item(81, 9).Value = 1200
$ws.Cells.Item(81, 11).Value = 3600
$ws.Cells.Item(81, 13).Value = -2477
$ws.Cells.Item(84, 8).Value = 27300
$ws.Cells.Item(84, 9).Value = 1200
$ws.Cells.Item(84, 11).Value = 10800
$ws.Cells.Item(84, 13).Value = -5184
$ws.Cells.Item(87, 8).Value = 9222.223
$ws.Cells.Item(87, 9).Value = 4600
$ws.Cells.Item(87, 11).Value = 13800
$ws.Cells.Item(87, 13).Value = -12552
$ws.Cells.Item(90, 8).Value = 9222.223
$ws.Cells.Item(90, 9).Value = 4600
$ws.Cells.Item(90, 11).Value = 41400
$ws.Cells.Item(90, 13).Value = -35160
$ws.Cells.Item(131, 8).Value = 1645.1082
$ws.Cells.Item(131, 9).Value = 2278.6
$ws.Cells.Item(131, 10).Value = 1213.1818
$ws.Cells.Item(131, 11).Value = 6835.799999999999
$ws.Cells.Item(131, 12).Value = 3639.5454
$ws.Cells.Item(131, 13).Value = -1795.799999999999
$ws.Cells.Item(131, 14).Value = -13719.5454

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 4427.8887
$ws.Cells.Item(126, 9).Value = 4206.364
$ws.Cells.Item(126, 10).Value = 4639.7827
$ws.Cells.Item(126, 11).Value = 12619.092
$ws.Cells.Item(126, 12).Value = 13919.3481
$ws.Cells.Item(126, 13).Value = -10149.092
$ws.Cells.Item(126, 14).Value = -18859.3481
$ws.Cells.Item(135, 8).Value = 27160
$ws.Cells.Item(135, 10).Value = 27160
$ws.Cells.Item(135, 12).Value = 27160
$ws.Cells.Item(135, 14).Value = -37300

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 3082.08
$ws.Cells.Item(132, 9).Value = 2364.889
$ws.Cells.Item(132, 10).Value = 3485.5
$ws.Cells.Item(132, 11).Value = 7094.667
$ws.Cells.Item(132, 12).Value = 10456.5
$ws.Cells.Item(132, 13).Value = -4564.667
$ws.Cells.Item(132, 14).Value = -15516.5
$ws.Cells.Item(136, 8).Value = 2148.5715
$ws.Cells.Item(136, 9).Value = 1571.7646
$ws.Cells.Item(136, 10).Value = 4600
$ws.Cells.Item(136, 11).Value = 4715.293799999999
$ws.Cells.Item(136, 12).Value = 13800
$ws.Cells.Item(136, 13).Value = -2165.293799999999
$ws.Cells.Item(136, 14).Value = -18900

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 2426.4375
$ws.Cells.Item(122, 9).Value = 2023.4
$ws.Cells.Item(122, 10).Value = 2782.0588
$ws.Cells.Item(122, 11).Value = 6070.200000000001
$ws.Cells.Item(122, 12).Value = 8346.1764
$ws.Cells.Item(122, 13).Value = -3620.200000000001
$ws.Cells.Item(122, 14).Value = -13246.1764
$ws.Cells.Item(126, 8).Value = 2214.5881
$ws.Cells.Item(126, 9).Value = 1957.6428
$ws.Cells.Item(126, 10).Value = 2394.45
$ws.Cells.Item(126, 11).Value = 5872.928400000001
$ws.Cells.Item(126, 12).Value = 7183.349999999999
$ws.Cells.Item(126, 13).Value = -3402.928400000001
$ws.Cells.Item(126, 14).Value = -12123.35
$ws.Cells.Item(132, 8).Value = 9694.706
$ws.Cells.Item(132, 9).Value = 3081.2
$ws.Cells.Item(132, 10).Value = 19142.572
$ws.Cells.Item(132, 11).Value = 9243.599999999999
$ws.Cells.Item(132, 12).Value = 57427.716
$ws.Cells.Item(132, 13).Value = -6713.599999999999
$ws.Cells.Item(132, 14).Value = -62487.716
$ws.Cells.Item(136, 8).Value = 2719.6667
$ws.Cells.Item(136, 9).Value = 1502.55
$ws.Cells.Item(136, 11).Value = 4507.65
$ws.Cells.Item(136, 13).Value = -1957.65
